# "Fruta / hortaliza, semanal"
# The weekly refresh re-ordered the Achicoria price rows (rows 2-21) of the
# Mercado Mayorista Lo Valledor de Santiago sheet. For every row the Fecha
# (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M), Origen (O) and
# Precio $/Kg (P) columns now hold a different record's data; the remaining
# columns (A,B,C,E,F,G,H,I,N,Q,R) are identical for every row and stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2;  D=44251; J=120; K=5000; L=5000; M=5000; O="Región Metropolitana"; P=312},
    @{Row=3;  D=44236; J=180; K=4000; L=4500; M=4167; O="Región Metropolitana"; P=260},
    @{Row=4;  D=44204; J=430; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=5;  D=44292; J=90;  K=6000; L=6000; M=6000; O="Región Metropolitana"; P=375},
    @{Row=6;  D=44846; J=250; K=5000; L=5000; M=5000; O="Provincia de Quillota"; P=312},
    @{Row=7;  D=44231; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=8;  D=44208; J=160; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=9;  D=44230; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=10; D=44210; J=340; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=11; D=44215; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=12; D=44882; J=70;  K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=13; D=44187; J=160; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=14; D=44188; J=210; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=15; D=44189; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=16; D=44873; J=250; K=8000; L=8000; M=8000; O="Provincia de Quillota"; P=500},
    @{Row=17; D=44883; J=180; K=7000; L=8000; M=7500; O="Provincia de Quillota"; P=469},
    @{Row=18; D=44232; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=19; D=44186; J=160; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=20; D=44855; J=70;  K=6000; L=7000; M=6500; O="Provincia de Quillota"; P=406},
    @{Row=21; D=44875; J=90;  K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
}
